# Apply updated attendance counts ("想去人数") and ticket-price status
# ("最低票价") cells to the "展览" (sheet 1) and "全部类型" (sheet 4) tabs,
# matching the scraped bilibili show data refresh.

$wb = $excel.ActiveWorkbook

# "展览" sheet (sheet1.xml)
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2").Value = 629
$ws1.Range("G2").Value = "不可售"
$ws1.Range("F3").Value = 1311
$ws1.Range("F4").Value = 1180
$ws1.Range("F5").Value = 14416
$ws1.Range("F6").Value = 16858
$ws1.Range("F8").Value = 126
$ws1.Range("F9").Value = 35
$ws1.Range("F10").Value = 50
$ws1.Range("F11").Value = 203
$ws1.Range("F12").Value = 26
$ws1.Range("F18").Value = 116
$ws1.Range("F20").Value = 1285
$ws1.Range("F23").Value = 50
$ws1.Range("F24").Value = 32
$ws1.Range("F25").Value = 4
$ws1.Range("F26").Value = 6911
$ws1.Range("F27").Value = 975
$ws1.Range("F28").Value = 27
$ws1.Range("F29").Value = 1144
$ws1.Range("F30").Value = 16
$ws1.Range("F32").Value = 5791
$ws1.Range("F33").Value = 119
$ws1.Range("F35").Value = 203
$ws1.Range("F36").Value = 4915

# "全部类型" sheet (sheet4.xml)
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F2").Value = 629
$ws4.Range("G2").Value = "不可售"
$ws4.Range("F3").Value = 1311
$ws4.Range("F4").Value = 1180
$ws4.Range("F5").Value = 14416
$ws4.Range("F6").Value = 16859
$ws4.Range("F8").Value = 126
$ws4.Range("F9").Value = 35
$ws4.Range("F10").Value = 50
$ws4.Range("F11").Value = 203
$ws4.Range("F12").Value = 26
$ws4.Range("F18").Value = 116
$ws4.Range("F20").Value = 1285
$ws4.Range("F24").Value = 50
$ws4.Range("F25").Value = 32
$ws4.Range("F26").Value = 4
$ws4.Range("F27").Value = 6911
$ws4.Range("F28").Value = 975
$ws4.Range("F29").Value = 27
$ws4.Range("F30").Value = 1144
$ws4.Range("F31").Value = 16
$ws4.Range("F34").Value = 5791
$ws4.Range("F35").Value = 119
$ws4.Range("F37").Value = 203
$ws4.Range("F38").Value = 4915

